$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "inclusao do google scholar": a new strategy row (E3) is inserted at row 4, which pushes
# the existing E4, E5, E6, E7 rows down by one (their own data stays attached to their labels).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the E3 strategy data.
$ws.Range("A4").Value = "E3"
$ws.Range("B4").Value = "E3"
$ws.Range("C4").Value = 0.02313624678663239
$ws.Range("D4").Value = 0.4090909090909091

# Make sure the new label cell (column A) keeps the same bold/centered/bordered look as the
# other strategy-label cells in column A.
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A4").VerticalAlignment = -4160    # xlTop
$ws.Range("A4").Borders.LineStyle = 1
